$d = $word.ActiveDocument

# Locate the full sentence run that ends in "...variables física."
$rng = $d.Content
$found = $rng.Find.Execute(
    "Se espera tener un manual de guía para aquellas personas o estudiantes que quieran inicializarse en el diseño de sistemas embebidos para ser aplicados a la internet de las cosas, o control de variables física.",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0
)

if (-not $found) {
    throw "Target sentence not found"
}

# Work off a fresh Range built from the found boundaries (InsertXML needs a
# Range that hasn't been left in "Find" mode to correctly replace content).
$target = $d.Range($rng.Start, $rng.End)

$openXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r w:rsidRPr="00173831"><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:eastAsia="es-ES"/></w:rPr><w:t>Se espera tener un manual de guía para aquellas personas o estudiantes que quieran inicializarse en el diseño de sistemas embebidos para ser aplicados a la internet de las cosas, o control de variables física</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:eastAsia="es-ES"/></w:rPr><w:t>s</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:eastAsia="es-ES"/></w:rPr><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($openXml)
